$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36/37 swap: NEARProtocol and FirstDigitalUSD change rank order ---
$ws.Range("B36").Value = "'FirstDigitalUSD"
$ws.Range("C36").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  +0.02%  "

$ws.Range("B37").Value = "'NEARProtocol"
$ws.Range("C37").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.82"
$ws.Range("E37").Value = "'  +6.53%  "

# --- Row 41/42 swap: Aave and USDe change rank order ---
$ws.Range("B41").Value = "'USDe"
$ws.Range("C41").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "'  -0.03%  "

$ws.Range("B42").Value = "'Aave"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'169.98"
$ws.Range("E42").Value = "'  +7.35%  "

# --- Remaining price/volume updates ---
$ws.Range("D2").Value = "'63.363.70"
$ws.Range("E2").Value = "'  +0.64%  "
$ws.Range("D3").Value = "'2.640.54"
$ws.Range("E3").Value = "'  +1.96%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'588.97"
$ws.Range("E5").Value = "'  +0.75%  "
$ws.Range("D6").Value = "'143.63"
$ws.Range("E6").Value = "'  -2.31%  "
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "'  -1.66%  "
$ws.Range("D9").Value = "'2.639.50"
$ws.Range("E9").Value = "'  +1.98%  "
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "'  -1.49%  "
$ws.Range("D11").Value = "'5.57"
$ws.Range("E11").Value = "'  -1.38%  "
$ws.Range("E12").Value = "'  +0.11%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "'  -0.52%  "
$ws.Range("D14").Value = "'27.31"
$ws.Range("E14").Value = "'  +0.10%  "
$ws.Range("D15").Value = "'3.115.33"
$ws.Range("E15").Value = "'  +2.01%  "
$ws.Range("D16").Value = "'63.340.53"
$ws.Range("E16").Value = "'  +0.74%  "
$ws.Range("E17").Value = "'  -1.54%  "
$ws.Range("D18").Value = "'2.644.73"
$ws.Range("E18").Value = "'  +2.27%  "
$ws.Range("D19").Value = "'11.30"
$ws.Range("E19").Value = "'  -0.03%  "
$ws.Range("D20").Value = "'339.99"
$ws.Range("E20").Value = "'  -0.58%  "
$ws.Range("D21").Value = "'4.34"
$ws.Range("E21").Value = "'  -1.25%  "
$ws.Range("D22").Value = "'6.67"
$ws.Range("E22").Value = "'  -0.24%  "
$ws.Range("E23").Value = "'  +0.04%  "
$ws.Range("D24").Value = "'67.46"
$ws.Range("E24").Value = "'  +0.31%  "
$ws.Range("D25").Value = "'1.64"
$ws.Range("E25").Value = "'  +3.01%  "
$ws.Range("D27").Value = "'0.165"
$ws.Range("E27").Value = "'  -1.28%  "
$ws.Range("D28").Value = "'544.84"
$ws.Range("E28").Value = "'  +14.48%  "
$ws.Range("E29").Value = "'  +0.15%  "
$ws.Range("D30").Value = "'8.38"
$ws.Range("E30").Value = "'  +0.42%  "
$ws.Range("D31").Value = "'7.73"
$ws.Range("E31").Value = "'  -1.48%  "
$ws.Range("E32").Value = "'  +13.01%  "
$ws.Range("E33").Value = "'  +1.71%  "
$ws.Range("D34").Value = "'0.0₃0803"
$ws.Range("E34").Value = "'  -2.40%  "
$ws.Range("D35").Value = "'173.42"
$ws.Range("E35").Value = "'  -1.95%  "
$ws.Range("D38").Value = "'0.400"
$ws.Range("E38").Value = "'  -1.24%  "
$ws.Range("D39").Value = "'19.00"
$ws.Range("E39").Value = "'  -0.21%  "
$ws.Range("E40").Value = "'  +3.61%  "
$ws.Range("D43").Value = "'40.22"
$ws.Range("E43").Value = "'  +1.91%  "
$ws.Range("D44").Value = "'3.71"
$ws.Range("E44").Value = "'  -1.25%  "
$ws.Range("D45").Value = "'22.13"
$ws.Range("E45").Value = "'  +3.74%  "
$ws.Range("D46").Value = "'0.626"
$ws.Range("E46").Value = "'  -1.24%  "
$ws.Range("D47").Value = "'0.0549"
$ws.Range("E47").Value = "'  +0.96%  "
$ws.Range("D48").Value = "'0.0956"
$ws.Range("E48").Value = "'  -1.47%  "
$ws.Range("E49").Value = "'  +0.08%  "
$ws.Range("D50").Value = "'18.63"
$ws.Range("E50").Value = "'  +1.61%  "
$ws.Range("E51").Value = "'  -0.80%  "

